$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.146479845046997
$ws.Range("B1").Value = 2.576543092727661
$ws.Range("C1").Value = 6.031661987304688
$ws.Range("D1").Value = 2.1337730884552
$ws.Range("E1").Value = 1.228684306144714
